$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 6; existing rows 6-48 shift down to 7-49
$ws.Rows.Item(6).Insert()

# Populate the newly inserted row 6 with the new weekly entry
$ws.Cells.Item(6, 1).Value = 10
$ws.Cells.Item(6, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(6, 3).Value = "La Araucanía"
$ws.Cells.Item(6, 4).Value = 44503
$ws.Cells.Item(6, 5).Value = 9
$ws.Cells.Item(6, 6).Value = 300000000
$ws.Cells.Item(6, 7).Value = "Espárragos"
$ws.Cells.Item(6, 8).Value = "Sin especificar"
$ws.Cells.Item(6, 9).Value = "Primera"
$ws.Cells.Item(6, 10).Value = 145
$ws.Cells.Item(6, 11).Value = 1200
$ws.Cells.Item(6, 12).Value = 1300
$ws.Cells.Item(6, 13).Value = 1245
$ws.Cells.Item(6, 14).Value = "$/kilo"
$ws.Cells.Item(6, 15).Value = "Región del Maule"
$ws.Cells.Item(6, 16).Value = 1245
$ws.Cells.Item(6, 17).Value = 1
$ws.Cells.Item(6, 18).Value = "Hortaliza"
